# Add a new worksheet named "Last" after the existing "Another data set" sheet,
# populate it with the example data, and make it the active sheet/selection
# (matching the "exporting only specific worksheets" example addition).

$wb = $excel.ActiveWorkbook

$lastExistingSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastExistingSheet)
$newSheet.Name = "Last"

# Header row
$newSheet.Range("A1").Value = "this"
$newSheet.Range("B1").Value = "is"
$newSheet.Range("C1").Value = "great"

# Data rows 2..30
for ($row = 2; $row -le 30; $row++) {
    $a = $row - 2
    $b = ($row - 1) * 10
    $c = 99 - ($row - 2) * 10
    $newSheet.Cells.Item($row, 1).Value = $a
    $newSheet.Cells.Item($row, 2).Value = $b
    $newSheet.Cells.Item($row, 3).Value = $c
}

# Make the new sheet the active one, with H18 selected.
$newSheet.Activate()
$newSheet.Range("H18").Select() | Out-Null
